$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 212016.6  # H86: was 133673.25
$ws.Cells.Item(86, 9).Value = 18859.666  # I86: was 12576.4
$ws.Cells.Item(86, 10).Value = 501752  # J86: was 335501.34
$ws.Cells.Item(86, 11).Value = 18859.666  # K86: was 12576.4
$ws.Cells.Item(86, 12).Value = 501752  # L86: was 335501.34
$ws.Cells.Item(86, 13).Value = -17736.666  # M86: was -11453.4
$ws.Cells.Item(86, 14).Value = -503998  # N86: was -337747.34
$ws.Cells.Item(89, 8).Value = 212016.6  # H89: was 133673.25
$ws.Cells.Item(89, 9).Value = 18859.666  # I89: was 12576.4
$ws.Cells.Item(89, 10).Value = 501752  # J89: was 335501.34
$ws.Cells.Item(89, 11).Value = 94298.33  # K89: was 62882
$ws.Cells.Item(89, 12).Value = 2508760  # L89: was 1677506.7
$ws.Cells.Item(89, 13).Value = -88682.33  # M89: was -57266
$ws.Cells.Item(89, 14).Value = -2519992  # N89: was -1688738.7
$ws.Cells.Item(103, 8).Value = 916.3333  # H103: was 1099.5
$ws.Cells.Item(103, 9).Value = 550  # I103: was 0
$ws.Cells.Item(103, 11).Value = 1650  # K103: was 0
$ws.Cells.Item(103, 13).Value = -1064  # M103: was None
$ws.Cells.Item(138, 8).Value = 4256.6772  # H138: was 4639.2593
$ws.Cells.Item(138, 10).Value = 5140.773  # J138: was 5911.1113
$ws.Cells.Item(138, 12).Value = 15422.319  # L138: was 17733.3339
$ws.Cells.Item(138, 14).Value = -25702.319  # N138: was -28013.3339

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 849.5  # H4: was 899.3333
$ws.Cells.Item(4, 9).Value = 700  # I4: was 849.5
$ws.Cells.Item(4, 11).Value = 700  # K4: was 849.5
$ws.Cells.Item(4, 13).Value = -584  # M4: was -733.5
$ws.Cells.Item(40, 8).Value = 35000  # H40: was 0
$ws.Cells.Item(40, 10).Value = 35000  # J40: was 0
$ws.Cells.Item(40, 12).Value = 35000  # L40: was 0
$ws.Cells.Item(40, 14).Value = -35352  # N40: was None
$ws.Cells.Item(61, 8).Value = 1337.3846  # H61: was 1383.8334
$ws.Cells.Item(61, 9).Value = 1337.3846  # I61: was 1383.8334
$ws.Cells.Item(61, 11).Value = 1337.3846  # K61: was 1383.8334
$ws.Cells.Item(61, 13).Value = -1125.3846  # M61: was -1171.8334
$ws.Cells.Item(136, 8).Value = 1337.3846  # H136: was 1383.8334
$ws.Cells.Item(136, 9).Value = 1337.3846  # I136: was 1383.8334
$ws.Cells.Item(136, 11).Value = 4012.1538  # K136: was 4151.5002
$ws.Cells.Item(136, 13).Value = -1462.1538  # M136: was -1601.5002

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(16, 8).Value = 760  # H16: was 1700
$ws.Cells.Item(16, 9).Value = 700  # I16: was 1000
$ws.Cells.Item(16, 10).Value = 850  # J16: was 2400
$ws.Cells.Item(16, 11).Value = 700  # K16: was 1000
$ws.Cells.Item(16, 12).Value = 850  # L16: was 2400
$ws.Cells.Item(16, 13).Value = -530  # M16: was -830
$ws.Cells.Item(16, 14).Value = -1190  # N16: was -2740
$ws.Cells.Item(86, 8).Value = 1511  # H86: was 1566.5
$ws.Cells.Item(86, 9).Value = 1480  # I86: was 1533.3334
$ws.Cells.Item(86, 11).Value = 1480  # K86: was 1533.3334
$ws.Cells.Item(86, 13).Value = -357  # M86: was -410.3334
$ws.Cells.Item(89, 8).Value = 1511  # H89: was 1566.5
$ws.Cells.Item(89, 9).Value = 1480  # I89: was 1533.3334
$ws.Cells.Item(89, 11).Value = 7400  # K89: was 7666.666999999999
$ws.Cells.Item(89, 13).Value = -1784  # M89: was -2050.666999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 3000  # H4: was 3833.3333
$ws.Cells.Item(4, 10).Value = 3000  # J4: was 8000
$ws.Cells.Item(4, 12).Value = 3000  # L4: was 8000
$ws.Cells.Item(4, 14).Value = -3224  # N4: was -8224
$ws.Cells.Item(12, 8).Value = 14993.333  # H12: was 11344.5
$ws.Cells.Item(12, 9).Value = 19992.5  # I12: was 13461
$ws.Cells.Item(12, 11).Value = 19992.5  # K12: was 13461
$ws.Cells.Item(12, 13).Value = -19822.5  # M12: was -13291
$ws.Cells.Item(14, 8).Value = 505  # H14: was 0
$ws.Cells.Item(14, 9).Value = 10  # I14: was 0
$ws.Cells.Item(14, 10).Value = 1000  # J14: was 0
$ws.Cells.Item(14, 11).Value = 10  # K14: was 0
$ws.Cells.Item(14, 12).Value = 1000  # L14: was 0
$ws.Cells.Item(14, 13).Value = 160  # M14: was None
$ws.Cells.Item(14, 14).Value = -1340  # N14: was None
$ws.Cells.Item(22, 8).Value = 398.66666  # H22: was 366
$ws.Cells.Item(22, 9).Value = 398.66666  # I22: was 366
$ws.Cells.Item(22, 11).Value = 398.66666  # K22: was 366
$ws.Cells.Item(22, 13).Value = -48.66665999999998  # M22: was -16
$ws.Cells.Item(33, 8).Value = 3625  # H33: was 1657.75
$ws.Cells.Item(33, 9).Value = 3400  # I33: was 1657.75
$ws.Cells.Item(33, 10).Value = 3962.5  # J33: was 0
$ws.Cells.Item(33, 11).Value = 3400  # K33: was 1657.75
$ws.Cells.Item(33, 12).Value = 3962.5  # L33: was 0
$ws.Cells.Item(33, 13).Value = -3021  # M33: was -1278.75
$ws.Cells.Item(33, 14).Value = -4720.5  # N33: was None
$ws.Cells.Item(42, 8).Value = 13499.5  # H42: was 0
$ws.Cells.Item(42, 9).Value = 9999  # I42: was 0
$ws.Cells.Item(42, 10).Value = 17000  # J42: was 0
$ws.Cells.Item(42, 11).Value = 9999  # K42: was 0
$ws.Cells.Item(42, 12).Value = 17000  # L42: was 0
$ws.Cells.Item(42, 13).Value = -9406  # M42: was None
$ws.Cells.Item(42, 14).Value = -18186  # N42: was None
$ws.Cells.Item(58, 8).Value = 2248  # H58: was 2198.4
$ws.Cells.Item(86, 8).Value = 14856.857  # H86: was 14592.5
$ws.Cells.Item(86, 9).Value = 13799.6  # I86: was 14948.2
$ws.Cells.Item(86, 10).Value = 17500  # J86: was 13999.667
$ws.Cells.Item(86, 11).Value = 13799.6  # K86: was 14948.2
$ws.Cells.Item(86, 12).Value = 17500  # L86: was 13999.667
$ws.Cells.Item(86, 13).Value = -12676.6  # M86: was -13825.2
$ws.Cells.Item(86, 14).Value = -19746  # N86: was -16245.667
$ws.Cells.Item(89, 8).Value = 14856.857  # H89: was 14592.5
$ws.Cells.Item(89, 9).Value = 13799.6  # I89: was 14948.2
$ws.Cells.Item(89, 10).Value = 17500  # J89: was 13999.667
$ws.Cells.Item(89, 11).Value = 68998  # K89: was 74741
$ws.Cells.Item(89, 12).Value = 87500  # L89: was 69998.33499999999
$ws.Cells.Item(89, 13).Value = -63382  # M89: was -69125
$ws.Cells.Item(89, 14).Value = -98732  # N89: was -81230.33499999999
$ws.Cells.Item(134, 8).Value = 2033.6207  # H134: was 2077.7144
$ws.Cells.Item(134, 9).Value = 2149.6086  # I134: was 2158.3044
$ws.Cells.Item(134, 10).Value = 1589  # J134: was 1707
$ws.Cells.Item(134, 11).Value = 6448.825800000001  # K134: was 6474.9132
$ws.Cells.Item(134, 12).Value = 4767  # L134: was 5121
$ws.Cells.Item(134, 13).Value = -3913.825800000001  # M134: was -3939.9132
$ws.Cells.Item(134, 14).Value = -9837  # N134: was -10191
$ws.Cells.Item(136, 8).Value = 2248  # H136: was 2198.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(29, 8).Value = 957  # H29: was 957.25
$ws.Cells.Item(29, 9).Value = 899  # I29: was 900
$ws.Cells.Item(29, 11).Value = 2697  # K29: was 2700
$ws.Cells.Item(29, 13).Value = -2420  # M29: was -2423
$ws.Cells.Item(34, 8).Value = 519.5  # H34: was 539
$ws.Cells.Item(34, 10).Value = 500  # J34: was 0
$ws.Cells.Item(34, 12).Value = 1500  # L34: was 0
$ws.Cells.Item(34, 14).Value = -1668  # N34: was None
$ws.Cells.Item(69, 8).Value = 999.6667  # H69: was 2000
$ws.Cells.Item(69, 9).Value = 999  # I69: was 0
$ws.Cells.Item(69, 10).Value = 1000  # J69: was 2000
$ws.Cells.Item(69, 11).Value = 2997  # K69: was 0
$ws.Cells.Item(69, 12).Value = 3000  # L69: was 6000
$ws.Cells.Item(69, 13).Value = -2186  # M69: was None
$ws.Cells.Item(69, 14).Value = -4622  # N69: was -7622
$ws.Cells.Item(72, 8).Value = 999.6667  # H72: was 2000
$ws.Cells.Item(72, 9).Value = 999  # I72: was 0
$ws.Cells.Item(72, 10).Value = 1000  # J72: was 2000
$ws.Cells.Item(72, 11).Value = 8991  # K72: was 0
$ws.Cells.Item(72, 12).Value = 9000  # L72: was 18000
$ws.Cells.Item(72, 13).Value = -4935  # M72: was None
$ws.Cells.Item(72, 14).Value = -17112  # N72: was -26112
$ws.Cells.Item(113, 8).Value = 1558.7  # H113: was 1507.909
$ws.Cells.Item(113, 9).Value = 1314.6666  # I113: was 1236
$ws.Cells.Item(113, 11).Value = 3943.9998  # K113: was 3708
$ws.Cells.Item(113, 13).Value = -1773.9998  # M113: was -1538
$ws.Cells.Item(131, 8).Value = 835571.5600000001  # H131: was 911484.4399999999
$ws.Cells.Item(131, 9).Value = 765  # I131: was 1000
$ws.Cells.Item(131, 11).Value = 2295  # K131: was 3000
$ws.Cells.Item(131, 13).Value = 2745  # M131: was 2040

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 7950  # H10: was 7100
$ws.Cells.Item(10, 9).Value = 9000  # I10: was 7500
$ws.Cells.Item(10, 11).Value = 9000  # K10: was 7500
$ws.Cells.Item(10, 13).Value = -8831  # M10: was -7331
$ws.Cells.Item(62, 8).Value = 0  # H62: was 90077
$ws.Cells.Item(62, 9).Value = 0  # I62: was 90077
$ws.Cells.Item(62, 11).Value = 0  # K62: was 90077
$ws.Cells.Item(62, 13).ClearContents()  # M62: was -89391
$ws.Cells.Item(65, 8).Value = 0  # H65: was 90077
$ws.Cells.Item(65, 9).Value = 0  # I65: was 90077
$ws.Cells.Item(65, 11).Value = 0  # K65: was 270231
$ws.Cells.Item(65, 13).ClearContents()  # M65: was -266799
$ws.Cells.Item(70, 8).Value = 5243.2  # H70: was 5304
$ws.Cells.Item(70, 10).Value = 5336.3335  # J70: was 5504.5
$ws.Cells.Item(70, 12).Value = 5336.3335  # L70: was 5504.5
$ws.Cells.Item(70, 14).Value = -5876.3335  # N70: was -6044.5
$ws.Cells.Item(73, 8).Value = 5243.2  # H73: was 5304
$ws.Cells.Item(73, 10).Value = 5336.3335  # J73: was 5504.5
$ws.Cells.Item(73, 12).Value = 5336.3335  # L73: was 5504.5
$ws.Cells.Item(73, 14).Value = -7208.3335  # N73: was -7376.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2361.182  # H46: was 2934.625
$ws.Cells.Item(46, 9).Value = 735  # I46: was 892.6667
$ws.Cells.Item(46, 10).Value = 3716.3333  # J46: was 4159.8
$ws.Cells.Item(46, 11).Value = 735  # K46: was 892.6667
$ws.Cells.Item(46, 12).Value = 3716.3333  # L46: was 4159.8
$ws.Cells.Item(46, 13).Value = -547  # M46: was -704.6667
$ws.Cells.Item(46, 14).Value = -4092.3333  # N46: was -4535.8
$ws.Cells.Item(74, 8).Value = 72500  # H74: was 83598.5
$ws.Cells.Item(74, 9).Value = 72500  # I74: was 83598.5
$ws.Cells.Item(74, 11).Value = 72500  # K74: was 83598.5
$ws.Cells.Item(74, 13).Value = -71502  # M74: was -82600.5
$ws.Cells.Item(77, 8).Value = 72500  # H77: was 83598.5
$ws.Cells.Item(77, 9).Value = 72500  # I77: was 83598.5
$ws.Cells.Item(77, 11).Value = 217500  # K77: was 250795.5
$ws.Cells.Item(77, 13).Value = -212508  # M77: was -245803.5
$ws.Cells.Item(100, 8).Value = 4262.875  # H100: was 4044.6667
$ws.Cells.Item(100, 10).Value = 6500  # J100: was 5099.6665
$ws.Cells.Item(100, 12).Value = 6500  # L100: was 5099.6665
$ws.Cells.Item(100, 14).Value = -7582  # N100: was -6181.6665

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 619999.7  # H2: was 99999
$ws.Cells.Item(2, 10).Value = 619999.7  # J2: was 99999
$ws.Cells.Item(2, 12).Value = 619999.7  # L2: was 99999
$ws.Cells.Item(2, 14).Value = -620223.7  # N2: was -100223
$ws.Cells.Item(70, 8).Value = 846575  # H70: was 84675
$ws.Cells.Item(70, 10).Value = 846575  # J70: was 84675
$ws.Cells.Item(70, 12).Value = 846575  # L70: was 84675
$ws.Cells.Item(70, 14).Value = -847205  # N70: was -85305
$ws.Cells.Item(73, 8).Value = 846575  # H73: was 84675
$ws.Cells.Item(73, 10).Value = 846575  # J73: was 84675
$ws.Cells.Item(73, 12).Value = 846575  # L73: was 84675
$ws.Cells.Item(73, 14).Value = -848759  # N73: was -86859
$ws.Cells.Item(126, 8).Value = 4480.909  # H126: was 4366
$ws.Cells.Item(126, 9).Value = 4481  # I126: was 4308.6
$ws.Cells.Item(126, 11).Value = 13443  # K126: was 12925.8
$ws.Cells.Item(126, 13).Value = -10973  # M126: was -10455.8
$ws.Cells.Item(136, 8).Value = 2599.125  # H136: was 2605.375
$ws.Cells.Item(136, 9).Value = 2932.3333  # I136: was 2940.6667
$ws.Cells.Item(136, 11).Value = 8796.999899999999  # K136: was 8822.000100000001
$ws.Cells.Item(136, 13).Value = -6246.999899999999  # M136: was -6272.000100000001
